$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")
$ws.Activate()
$win = $excel.ActiveWindow
$sv = $win.ActiveSheetView
$props = $sv | Get-Member
Write-Output $props
